# Update likelihood formula results (columns D:L) for rows 2-16
# to reflect a re-run of the likelihood computation with new random
# train/test splits (plots methods and prints).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1.027, 5.3, 3.183333333333334, 1.338055555555555, 1.266146835245512, 0.216697217897322, 0.06465291590195688, 0.7510578991915623, 0.02427904159696624),
    @(1.239, 4.0, 3.183333333333334, 1.338055555555555, 1.438804582540721, 0.216697217897322, 0.2688050773771029, 0.78158932586104, 0.1050475896075973),
    @(1.166, 2.5, 3.183333333333334, 1.338055555555555, 1.638025060189038, 0.216697217897322, 0.2896640519617951, 0.5125226083978361, 0.07422968773527278),
    @(0.9089999999999999, 3.0, 3.183333333333334, 1.338055555555555, 1.571618234306266, 0.216697217897322, 0.3405794120351475, 0.3111778354464314, 0.05299038211735774),
    @(1.1865, 2.5, 3.183333333333334, 1.338055555555555, 1.638025060189038, 0.216697217897322, 0.2896640519617951, 0.5354083298243082, 0.07754427313550316),
    @(1.89, 1.8, 3.183333333333334, 1.338055555555555, 1.73099461642492, 0.216697217897322, 0.1687024282500931, 0.8084407205555387, 0.06819295632698717),
    @(2.171, 5.0, 3.638888888888889, 1.517098765432099, 1.305990930775175, 0.216697217897322, 0.1758859568782206, 0.1524715479054494, 0.01340880205002671),
    @(1.8415, 2.9, 3.638888888888889, 1.517098765432099, 1.58489959948282, 0.216697217897322, 0.2705568796570144, 0.7362120005072895, 0.09959361081165026),
    @(1.031, 2.95, 3.638888888888889, 1.517098765432099, 1.578258916894543, 0.216697217897322, 0.2769980447648336, 0.4294069895091945, 0.05947244825120014),
    @(1.742, 3.0, 3.638888888888889, 1.517098765432099, 1.571618234306266, 0.216697217897322, 0.2831256128309231, 0.8014808763750364, 0.1134598821479738),
    @(2.302, 1.7, 3.638888888888889, 1.517098765432099, 1.744275981601474, 0.216697217897322, 0.0938261323161729, 0.4181011246514542, 0.01961440572154402),
    @(1.2415, 4.9, 3.638888888888889, 1.517098765432099, 1.31927229595173, 0.216697217897322, 0.1917625863465881, 0.8451276423576342, 0.08103193124574712),
    @(1.07, 4.25, 3.638888888888889, 1.517098765432099, 1.405601169599334, 0.216697217897322, 0.2863842059988284, 0.6608780899792053, 0.09463252353035848),
    @(2.446, 2.55, 3.638888888888889, 1.517098765432099, 1.631384377600761, 0.216697217897322, 0.2191255786365801, 0.1853560074436751, 0.02030812119243078),
    @(1.402, 5.5, 3.638888888888889, 1.517098765432099, 1.239584104892403, 0.216697217897322, 0.1034255012828138, 0.8063984660251935, 0.04170108279117386)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $ws.Cells.Item($r, 4 + $j).Value = $vals[$j]
    }
}

Write-Output "Updated likelihood columns D:L for rows 2-16"